# Eliminan EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Periodo Mora" column (E) values for rows 16-19 are reordered (reversed)
$ws.Range("E16").Value = "2408"
$ws.Range("E17").Value = "2407"
$ws.Range("E18").Value = "2406"
$ws.Range("E19").Value = "2402"

# "Valor Mora" column (F) values for rows 16 and 19 are swapped accordingly
$ws.Range("F16").Value = 42000
$ws.Range("F19").Value = 29466
